$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing columns B,C,D -> C,D,E
$ws.Range("A1").EntireColumn.Insert()

# New B1 header "segments" - copy format (border/bold/alignment) from C1 (old B1)
$ws.Range("C1").Copy() | Out-Null
$ws.Range("B1").PasteSpecial(-4122) | Out-Null
$ws.Range("B1").Value = "segments"

# Fill new column A (the 0-based segment index) for rows 2..20,
# copying the "index" style (border/bold/alignment) from column B (old column A)
for ($i = 2; $i -le 20; $i++) {
    $ws.Range("B$i").Copy() | Out-Null
    $ws.Range("A$i").PasteSpecial(-4122) | Out-Null
    $ws.Range("A$i").Value = $i - 2
}

# The segment-name column (now B) should no longer carry the bordered/bold style
$ws.Range("B2:B20").ClearFormats()

$excel.CutCopyMode = $false
